# Apply the "SortareRezultateForms/raspunsuri.xlsx" revision:
#  - refresh the "Date studenti" (previous-year averages) lookup table with
#    a new data pull (new order + new values per student)
#  - drop the now-unused "Discipline" helper sheet
#  - leave "Raspunsuri la formular" as the active/selected sheet, and park
#    the selection on "Repartizare" at a single cell instead of the old
#    full-table selection

$wb = $excel.ActiveWorkbook

$wsRaspunsuri   = $wb.Worksheets.Item("Raspunsuri la formular")
$wsDateStudenti = $wb.Worksheets.Item("Date studenti")
$wsRepartizare  = $wb.Worksheets.Item("Repartizare")
$wsDiscipline   = $wb.Worksheets.Item("Discipline")

# --- "Date studenti": new pull of last-year averages, rows 2-21 ------------
$studentData = @(
    @("antonio.popescu02@e-uvt.ro", 8.691),
    @("darius.hoalba02@e-uvt.ro", 8.251),
    @("cezar.petreanu02@e-uvt.ro", 8.314),
    @("grigorie.smarandache02@e-uvt.ro", 9.219),
    @("robert.ionita02@e-uvt.ro", 6.981),
    @("stefan.soare02@e-uvt.ro", 9.069),
    @("alexandru.bran02@e-uvt.ro", 7.407),
    @("octavian.ilies02@e-uvt.ro", 8.888),
    @("damian.cernea00@e-uvt.ro", 7.747),
    @("costel.anghel02@e-uvt.ro", 6.382),
    @("vlad.gozman02@e-uvt.ro", 5.675),
    @("mihai.tudor02@e-uvt.ro", 5.018),
    @("vlad.bradea02@e-uvt.ro", 5.796),
    @("catalin.baltaretu03@e-uvt.ro", 8.88),
    @("danut.termure03@e-uvt.ro", 7.167),
    @("petre.raduletu02@e-uvt.ro", 8.99),
    @("cristian.gusatu02@e-uvt.ro", 7.624),
    @("alex.falcuta02@e-uvt.ro", 8.107),
    @("tudor.mateian02@e-uvt.ro", 5.4),
    @("andreea.brad02@e-uvt.ro", 9.889)
)

$r = 2
foreach ($row in $studentData) {
    $wsDateStudenti.Cells.Item($r, 1).Value = $row[0]
    $wsDateStudenti.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# --- drop the "Discipline" helper sheet -------------------------------------
$wsDiscipline.Delete()

# --- selections / active sheet ---------------------------------------------
$wsRepartizare.Range("E16").Select()
$wsRaspunsuri.Activate()
$wsRaspunsuri.Range("A1:B1").Select()
